$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.4115256666666666
$ws.Cells.Item(2, 8).Value = 1.234577
$ws.Cells.Item(2, 9).Value = 0.2245998342667577
$ws.Cells.Item(2, 10).Value = 0.2245998342667577
$ws.Cells.Item(2, 13).Value = 2.883158333333334
$ws.Cells.Item(2, 14).Value = 8.649475
$ws.Cells.Item(2, 15).Value = 0.1005826776766536
$ws.Cells.Item(2, 16).Value = 0.1005826776766536
$ws.Cells.Item(2, 17).Value = 1.186493655230556
$ws.Cells.Item(2, 18).Value = 10.678442897075
$ws.Cells.Item(2, 19).Value = 0.0225908527362831
$ws.Cells.Item(2, 20).Value = 0.0225908527362831
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.4115256666666666
$ws.Cells.Item(3, 8).Value = 1.234577
$ws.Cells.Item(3, 9).Value = 0.2245998342667577
$ws.Cells.Item(3, 10).Value = 0.2245998342667577
$ws.Cells.Item(3, 15).Value = 0.4854237085598054
$ws.Cells.Item(3, 16).Value = 0.4854237085598054
$ws.Cells.Item(3, 17).Value = 5.726156467580111
$ws.Cells.Item(3, 18).Value = 51.535408208221
$ws.Cells.Item(3, 19).Value = 0.1090260844916872
$ws.Cells.Item(3, 20).Value = 0.1090260844916872
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.4115256666666666
$ws.Cells.Item(4, 8).Value = 1.234577
$ws.Cells.Item(4, 9).Value = 0.2245998342667577
$ws.Cells.Item(4, 10).Value = 0.2245998342667577
$ws.Cells.Item(4, 13).Value = 11.86694533333333
$ws.Cells.Item(4, 14).Value = 35.600836
$ws.Cells.Item(4, 15).Value = 0.413993613763541
$ws.Cells.Item(4, 16).Value = 0.413993613763541
$ws.Cells.Item(4, 17).Value = 4.883552589596889
$ws.Cells.Item(4, 18).Value = 43.951973306372
$ws.Cells.Item(4, 19).Value = 0.09298289703878743
$ws.Cells.Item(4, 20).Value = 0.09298289703878743
$ws.Cells.Item(5, 9).Value = 0.3944722233087159
$ws.Cells.Item(5, 10).Value = 0.3944722233087159
$ws.Cells.Item(5, 13).Value = 2.883158333333334
$ws.Cells.Item(5, 14).Value = 8.649475
$ws.Cells.Item(5, 15).Value = 0.1005826776766536
$ws.Cells.Item(5, 16).Value = 0.1005826776766536
$ws.Cells.Item(5, 17).Value = 2.083878608586112
$ws.Cells.Item(5, 18).Value = 18.754907477275
$ws.Cells.Item(5, 19).Value = 0.03967707248945348
$ws.Cells.Item(5, 20).Value = 0.03967707248945348
$ws.Cells.Item(6, 9).Value = 0.3944722233087159
$ws.Cells.Item(6, 10).Value = 0.3944722233087159
$ws.Cells.Item(6, 15).Value = 0.4854237085598054
$ws.Cells.Item(6, 16).Value = 0.4854237085598054
$ws.Cells.Item(6, 19).Value = 0.1914861695623486
$ws.Cells.Item(6, 20).Value = 0.1914861695623486
$ws.Cells.Item(7, 9).Value = 0.3944722233087159
$ws.Cells.Item(7, 10).Value = 0.3944722233087159
$ws.Cells.Item(7, 13).Value = 11.86694533333333
$ws.Cells.Item(7, 14).Value = 35.600836
$ws.Cells.Item(7, 15).Value = 0.413993613763541
$ws.Cells.Item(7, 16).Value = 0.413993613763541
$ws.Cells.Item(7, 17).Value = 8.577147235893777
$ws.Cells.Item(7, 18).Value = 77.194325123044
$ws.Cells.Item(7, 19).Value = 0.1633089812569138
$ws.Cells.Item(7, 20).Value = 0.1633089812569138
$ws.Cells.Item(8, 7).Value = 0.6979596666666668
$ws.Cells.Item(8, 8).Value = 2.093879
$ws.Cells.Item(8, 9).Value = 0.3809279424245264
$ws.Cells.Item(8, 10).Value = 0.3809279424245264
$ws.Cells.Item(8, 13).Value = 2.883158333333334
$ws.Cells.Item(8, 14).Value = 8.649475
$ws.Cells.Item(8, 15).Value = 0.1005826776766536
$ws.Cells.Item(8, 16).Value = 0.1005826776766536
$ws.Cells.Item(8, 17).Value = 2.012328229280556
$ws.Cells.Item(8, 18).Value = 18.110954063525
$ws.Cells.Item(8, 19).Value = 0.03831475245091698
$ws.Cells.Item(8, 20).Value = 0.03831475245091698
$ws.Cells.Item(9, 7).Value = 0.6979596666666668
$ws.Cells.Item(9, 8).Value = 2.093879
$ws.Cells.Item(9, 9).Value = 0.3809279424245264
$ws.Cells.Item(9, 10).Value = 0.3809279424245264
$ws.Cells.Item(9, 15).Value = 0.4854237085598054
$ws.Cells.Item(9, 16).Value = 0.4854237085598054
$ws.Cells.Item(9, 17).Value = 9.711730234874114
$ws.Cells.Item(9, 18).Value = 87.40557211386702
$ws.Cells.Item(9, 19).Value = 0.1849114545057696
$ws.Cells.Item(9, 20).Value = 0.1849114545057696
$ws.Cells.Item(10, 7).Value = 0.6979596666666668
$ws.Cells.Item(10, 8).Value = 2.093879
$ws.Cells.Item(10, 9).Value = 0.3809279424245264
$ws.Cells.Item(10, 10).Value = 0.3809279424245264
$ws.Cells.Item(10, 13).Value = 11.86694533333333
$ws.Cells.Item(10, 14).Value = 35.600836
$ws.Cells.Item(10, 15).Value = 0.413993613763541
$ws.Cells.Item(10, 16).Value = 0.413993613763541
$ws.Cells.Item(10, 17).Value = 8.282649209204891
$ws.Cells.Item(10, 18).Value = 74.54384288284402
$ws.Cells.Item(10, 19).Value = 0.1577017354678398
$ws.Cells.Item(10, 20).Value = 0.1577017354678398
